# Generate Report for Handoff
# Applies updated handoff/handback status + reordered rows to the
# localization-status workbook (Overview, zh-cn, de-de sheets).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("A2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.md'
$wsOverview.Range("B2").Value = 'e2e\cb49395a-ed34-4f50-82b2-7615fe5cb702.md'
$wsOverview.Range("G2").Value = '2016-09-02 10:28:38'
$wsOverview.Range("A3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.md'
$wsOverview.Range("B3").Value = 'e2e\e981ec47-db38-4d2c-b209-011a1b3e5640.md'
$wsOverview.Range("G3").Value = '2016-09-02 10:28:38'
$wsOverview.Range("A4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md'
$wsOverview.Range("B4").Value = 'e2e\008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md'
$wsOverview.Range("E4").Value = 'Ready for handoff'
$wsOverview.Range("F4").Value = 'Ready for handoff'
$wsOverview.Range("G4").Value = '2016-09-02 10:31:43'
$wsOverview.Range("A5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
$wsOverview.Range("B5").Value = 'e2e\c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
$wsOverview.Range("E5").Value = 'Ready for handoff'
$wsOverview.Range("F5").Value = 'Ready for handoff'
$wsOverview.Range("G5").Value = '2016-09-02 10:31:43'

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("A2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.md'
$wsZhCn.Range("G2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.982b4f776e430e268698c5db0af0a9adbcc9ef76.zh-cn.xlf'
$wsZhCn.Range("H2").Value = '2016-09-02 10:28:33'
$wsZhCn.Range("I2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.md'
$wsZhCn.Range("J2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.982b4f776e430e268698c5db0af0a9adbcc9ef76.zh-cn.xlf'
$wsZhCn.Range("K2").Value = '2016-09-02 10:29:43'
$wsZhCn.Range("A3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.md'
$wsZhCn.Range("G3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.bdf0ab7b67e639b7cb1d29f21124a864d97265e0.zh-cn.xlf'
$wsZhCn.Range("H3").Value = '2016-09-02 10:28:33'
$wsZhCn.Range("I3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.md'
$wsZhCn.Range("J3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.bdf0ab7b67e639b7cb1d29f21124a864d97265e0.zh-cn.xlf'
$wsZhCn.Range("K3").Value = '2016-09-02 10:29:43'
$wsZhCn.Range("A4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md'
$wsZhCn.Range("C4").Value = 'Ready for handoff'
$wsZhCn.Range("G4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.59f1e7bcba884cc0962bea7e951600d15ecfb5d2.zh-cn.xlf'
$wsZhCn.Range("H4").Value = '2016-09-02 10:31:38'
$wsZhCn.Range("I4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md'
$wsZhCn.Range("J4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.59f1e7bcba884cc0962bea7e951600d15ecfb5d2.zh-cn.xlf'
$wsZhCn.Range("K4").Value = '2016-09-02 10:31:00'
$wsZhCn.Range("P4").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f18ff73b15f3f7fd57a6741f248ad141172740f9/e2e/008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2654f8ba004cdf8da46bb2bc86ea9a5e30194cd/e2e/008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md.'
$wsZhCn.Range("A5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
$wsZhCn.Range("C5").Value = 'Ready for handoff'
$wsZhCn.Range("G5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.8204a07c236d29fb3a3d00e4b9372697cf33da29.zh-cn.xlf'
$wsZhCn.Range("H5").Value = '2016-09-02 10:31:38'
$wsZhCn.Range("I5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
$wsZhCn.Range("J5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.8204a07c236d29fb3a3d00e4b9372697cf33da29.zh-cn.xlf'
$wsZhCn.Range("K5").Value = '2016-09-02 10:31:00'
$wsZhCn.Range("P5").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f18ff73b15f3f7fd57a6741f248ad141172740f9/e2e/c680de99-2b03-4147-b4a1-f8eee182ff2d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2654f8ba004cdf8da46bb2bc86ea9a5e30194cd/e2e/c680de99-2b03-4147-b4a1-f8eee182ff2d.md.'

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("A2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.md'
$wsDeDe.Range("G2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.982b4f776e430e268698c5db0af0a9adbcc9ef76.de-de.xlf'
$wsDeDe.Range("H2").Value = '2016-09-02 10:28:38'
$wsDeDe.Range("I2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.md'
$wsDeDe.Range("J2").Value = 'cb49395a-ed34-4f50-82b2-7615fe5cb702.982b4f776e430e268698c5db0af0a9adbcc9ef76.de-de.xlf'
$wsDeDe.Range("K2").Value = '2016-09-02 10:29:50'
$wsDeDe.Range("A3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.md'
$wsDeDe.Range("G3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.bdf0ab7b67e639b7cb1d29f21124a864d97265e0.de-de.xlf'
$wsDeDe.Range("H3").Value = '2016-09-02 10:28:38'
$wsDeDe.Range("I3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.md'
$wsDeDe.Range("J3").Value = 'e981ec47-db38-4d2c-b209-011a1b3e5640.bdf0ab7b67e639b7cb1d29f21124a864d97265e0.de-de.xlf'
$wsDeDe.Range("K3").Value = '2016-09-02 10:29:50'
$wsDeDe.Range("A4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md'
$wsDeDe.Range("C4").Value = 'Ready for handoff'
$wsDeDe.Range("G4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.59f1e7bcba884cc0962bea7e951600d15ecfb5d2.de-de.xlf'
$wsDeDe.Range("H4").Value = '2016-09-02 10:31:43'
$wsDeDe.Range("I4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md'
$wsDeDe.Range("J4").Value = '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.59f1e7bcba884cc0962bea7e951600d15ecfb5d2.de-de.xlf'
$wsDeDe.Range("K4").Value = '2016-09-02 10:31:16'
$wsDeDe.Range("P4").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f18ff73b15f3f7fd57a6741f248ad141172740f9/e2e/008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2654f8ba004cdf8da46bb2bc86ea9a5e30194cd/e2e/008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md.'
$wsDeDe.Range("A5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
$wsDeDe.Range("C5").Value = 'Ready for handoff'
$wsDeDe.Range("G5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.8204a07c236d29fb3a3d00e4b9372697cf33da29.de-de.xlf'
$wsDeDe.Range("H5").Value = '2016-09-02 10:31:43'
$wsDeDe.Range("I5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
$wsDeDe.Range("J5").Value = 'c680de99-2b03-4147-b4a1-f8eee182ff2d.8204a07c236d29fb3a3d00e4b9372697cf33da29.de-de.xlf'
$wsDeDe.Range("K5").Value = '2016-09-02 10:31:16'
$wsDeDe.Range("P5").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f18ff73b15f3f7fd57a6741f248ad141172740f9/e2e/c680de99-2b03-4147-b4a1-f8eee182ff2d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2654f8ba004cdf8da46bb2bc86ea9a5e30194cd/e2e/c680de99-2b03-4147-b4a1-f8eee182ff2d.md.'

# --- Error Detail column widened to fit the new long diagnostic text --
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14

# --- Hyperlink display text follows the new row order (Address/target is
#     untouched - only the cached "display" text is refreshed), matching
#     the cyclic reshuffle of the 4 source files on every sheet. ---------
$overviewLinkText = @(
    'e2e\cb49395a-ed34-4f50-82b2-7615fe5cb702.md',
    'e2e\e981ec47-db38-4d2c-b209-011a1b3e5640.md',
    'e2e\008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md',
    'e2e\c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
)
$i = 0
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $overviewLinkText[$i]
    $i = $i + 1
}

$langLinkText = @(
    'cb49395a-ed34-4f50-82b2-7615fe5cb702.md',
    'cb49395a-ed34-4f50-82b2-7615fe5cb702.md',
    'e981ec47-db38-4d2c-b209-011a1b3e5640.md',
    'e981ec47-db38-4d2c-b209-011a1b3e5640.md',
    '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md',
    '008b9c1a-1e33-4258-9bdb-ea1dec4c3e6a.md',
    'c680de99-2b03-4147-b4a1-f8eee182ff2d.md',
    'c680de99-2b03-4147-b4a1-f8eee182ff2d.md'
)

$i = 0
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $langLinkText[$i]
    $i = $i + 1
}

$i = 0
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $langLinkText[$i]
    $i = $i + 1
}
